$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.564.47"
$ws.Range("E2").Value = "  +0.16%  "

# Row 3
$ws.Range("D3").Value = "3.047.07"
$ws.Range("E3").Value = "  -0.64%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.82"
$ws.Range("E5").Value = "  +0.34%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.73"
$ws.Range("E6").Value = "  -0.33%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.10%  "

# Row 8
$ws.Range("D8").Value = "3.044.75"
$ws.Range("E8").Value = "  -0.62%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.502"
$ws.Range("E9").Value = "  -0.03%  "

# Row 10
$ws.Range("E10").Value = "  +0.90%  "

# Row 11
$ws.Range("E11").Value = "  -6.82%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.467"
$ws.Range("E12").Value = "  +2.17%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000227"
$ws.Range("E13").Value = "  -0.45%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.55"
$ws.Range("E14").Value = "  -0.79%  "

# Row 15
$ws.Range("D15").Value = "3.548.76"
$ws.Range("E15").Value = "  -0.52%  "

# Row 16
$ws.Range("D16").Value = "63.586.90"
$ws.Range("E16").Value = "  +0.19%  "

# Row 17
$ws.Range("D17").Value = "3.049.63"
$ws.Range("E17").Value = "  -0.72%  "

# Row 18
$ws.Range("E18").Value = "  +0.50%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.66"
$ws.Range("E19").Value = "  -1.80%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "475.79"
$ws.Range("E20").Value = "  -1.55%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.89"
$ws.Range("E21").Value = "  +0.03%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.670"
$ws.Range("E22").Value = "  -1.02%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.47"
$ws.Range("E23").Value = "  +2.61%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.99"
$ws.Range("E24").Value = "  +8.98%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.76"
$ws.Range("E25").Value = "  -0.17%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.07%  "

# Row 27
$ws.Range("E27").Value = "  -0.62%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.89"
$ws.Range("E28").Value = "  -0.45%  "

# Row 29
$ws.Range("E29").Value = "  +0.67%  "

# Row 30
$ws.Range("E30").Value = "  -0.05%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.02"
$ws.Range("E31").Value = "  -0.49%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.14"
$ws.Range("E32").Value = "  -2.43%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.44"
$ws.Range("E33").Value = "  -0.38%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.55"
$ws.Range("E34").Value = "  -2.56%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.12"
$ws.Range("E35").Value = "  +1.92%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.96"
$ws.Range("E36").Value = "  -0.73%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0403"
$ws.Range("E37").Value = "  +1.59%  "

# Row 38
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "437.54"
$ws.Range("E38").Value = "  -6.33%  "

# Row 39
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.89"
$ws.Range("E39").Value = "  +12.40%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0800"
$ws.Range("E40").Value = "  -2.83%  "

# Row 41
$ws.Range("D41").Value = "2.940.97"
$ws.Range("E41").Value = "  -2.51%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.13"
$ws.Range("E42").Value = "  -0.94%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.113"
$ws.Range("E43").Value = "  -5.23%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.95"
$ws.Range("E44").Value = "  +0.95%  "

# Row 45
$ws.Range("E45").Value = "  +0.20%  "

# Row 47
$ws.Range("E47").Value = "  +2.52%  "

# Row 48
$ws.Range("E48").Value = "  +0.74%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "116.59"
$ws.Range("E49").Value = "  +0.00%  "

# Row 50
$ws.Range("D50").Value = "0.0₃0509"
$ws.Range("E50").Value = "  +0.12%  "

# Row 51
$ws.Range("E51").Value = "  -1.41%  "
